# country_comparison.xlsx update
# - Reorders the "country" columns (D1:N1) behind the header row (the
#   underlying shared-string table was re-sorted, which changes which
#   country label each column of data corresponds to).
# - Refreshes rows 2-6 with new percentages (recomputed without Russia,
#   i.e. "remove fast RU"), clearing a couple of cells that no longer have
#   data.
# - Appends three brand-new survey-question rows (7, 8, 9) with their own
#   data, some cells left blank where no data exists.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row: columns D1:N1 now point at a reshuffled set of countries.
# ---------------------------------------------------------------------
$headerCols = @("D","E","F","G","H","I","J","K","L","M","N")
$headerVals = @("France","Germany","Italy","Poland","Spain","United Kingdom","Switzerland","Japan","Russia","Saudi Arabia","USA")

for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + "1").Value = $headerVals[$i]
}

# ---------------------------------------------------------------------
# Row labels (column A) for the data rows.
# ---------------------------------------------------------------------
$labelRow2 = @'
Supports tax on world top 1% to finance global poverty reduction
(Additional 15% tax on income over [$120k/year in PPP])
'@

$labelRow3 = @'
Supports tax on world top 3% to finance global poverty reduction
(Additional 15% tax over [$80k], 30% over [$120k], 45% over [$1M])
'@

$labelRow4 = 'Prefers sustainable future'

$labelRow5 = @'
"Governments should actively cooperate to have all countries
converge in terms of GDP per capita by the end of the century"
'@

$labelRow6 = @'
Would support a global movement to tackle CC, tax millionaires,
 and fund LICs (either petition, demonstrate, strike, or donate)
'@

$labelRow7 = @'
More likely to vote for party if part of worldwide
coalition for climate action and global redistribution
'@

$labelRow8 = @'
Supports reparations for colonization and slavery in
the form of funding education and technology transfers
'@

$labelRow9 = '"My taxes should go towards solving global problems"'

$ws.Range("A2").Value = $labelRow2
$ws.Range("A3").Value = $labelRow3
$ws.Range("A4").Value = $labelRow4
$ws.Range("A5").Value = $labelRow5
$ws.Range("A6").Value = $labelRow6
$ws.Range("A7").Value = $labelRow7
$ws.Range("A8").Value = $labelRow8
$ws.Range("A9").Value = $labelRow9

# ---------------------------------------------------------------------
# Data values for columns B:N, rows 2-9. $null entries become blank
# cells (no cached value / type), matching cells that now have no data.
# ---------------------------------------------------------------------
$dataCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N")

$rowData = @{
    2 = @(0.241205184312405, 0.237490071485306, 0.243781094527363, 0.242201834862385, 0.145118733509235, 0.219512195121951, 0.218354430379747, 0.277777777777778, 0.341880341880342, 0.2, 0.204545454545455, 0.157446808510638, 0.315463917525773)
    3 = @(0.281590234246123, 0.295729250604351, 0.252525252525253, 0.312127236580517, 0.254641909814324, 0.240157480314961, 0.289198606271777, 0.267441860465116, 0.519148936170213, 0.284466019417476, 0.197894736842105, 0.145283018867925, 0.329449838187702)
    4 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    5 = @(0.279180423383181, 0.22951582867784, 0.231611893583725, 0.246376811594203, 0.12778603268945, 0.156043956043956, 0.158790170132325, 0.340974212034384, 0.340740740740741, 0.298397040690506, 0.231768231768232, 0.0743801652892562, 0.438388625592417)
    6 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, $null, 0, 0)
    7 = @(0.1704, 0.16, 0.171679197994987, 0.157442748091603, 0.111111111111111, 0.166, 0.129353233830846, 0.174334140435835, 0.232409381663113, 0.1765, $null, $null, 0.183666666666667)
    8 = @(0.414450291565922, 0.384271892830563, 0.422305764411028, 0.425572519083969, 0.247354497354497, $null, 0.386401326699834, 0.418886198547215, $null, $null, $null, $null, 0.455)
    9 = @(0.268977585201233, 0.2838, 0.422305764411028, 0.270038167938931, 0.16005291005291, 0.244, 0.203980099502488, 0.328087167070218, 0.345415778251599, 0.227, 0.297702297702298, 0.09, 0.322333333333333)
}

foreach ($r in 2..9) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $dataCols.Length; $i++) {
        $cellRef = $dataCols[$i] + $r
        if ($vals[$i] -eq $null) {
            $ws.Range($cellRef).ClearContents()
        } else {
            $ws.Range($cellRef).Value = $vals[$i]
        }
    }
}
